# Apply the weekly FlashScore data refresh:
#  1) Remove the QATAR - QSL match (old row 6: Al-Sadd vs Al Rayyan). All the
#     following rows (old rows 7-9) shift up by one automatically, becoming
#     the new rows 6-8.
#  2) Update the odds values that changed for every remaining match
#     (rows 2, 4, 5, 6, 7 and 8 after the shift) to reflect the latest
#     odds pulled from FlashScore.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Delete the Qatar QSL match row -----------------------------------
$ws.Rows("6").Delete()

# --- 2) Row 2 (Genoa - Fiorentina) ---------------------------------------
$ws.Range("Q2").Value = 2.04
$ws.Range("R2").Value = 1.86
$ws.Range("AE2").Value = 17

# --- 3) Row 4 (Como - Lazio) ----------------------------------------------
$ws.Range("G4").Value = 3.3
$ws.Range("H4").Value = 3.2
$ws.Range("I4").Value = 2.3
$ws.Range("J4").Value = 3.75
$ws.Range("L4").Value = 3
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("Y4").Value = 12
$ws.Range("AA4").Value = 26
$ws.Range("AD4").Value = 6.5
$ws.Range("AE4").Value = 15
$ws.Range("AH4").Value = 7.5
$ws.Range("AK4").Value = 21
$ws.Range("BC4").Value = 126

# --- 4) Row 5 (Modern Sport - Enppi) --------------------------------------
$ws.Range("H5").Value = 2.47
$ws.Range("L5").Value = 4.45
$ws.Range("M5").Value = 1.2
$ws.Range("O5").Value = 1.8
$ws.Range("R5").Value = 1.29
$ws.Range("AA5").Value = 30
$ws.Range("AB5").Value = 60
$ws.Range("AH5").Value = 6.3
$ws.Range("AJ5").Value = 14.5
$ws.Range("AL5").Value = 55
$ws.Range("AM5").Value = 90
$ws.Range("AP5").Value = 30
$ws.Range("AQ5").Value = 75
$ws.Range("AR5").Value = 150
$ws.Range("AZ5").Value = 175
$ws.Range("BA5").Value = 300

# --- 5) Row 6 after the shift (Al Ittihad - Al Ahli SC) ------------------
$ws.Range("G6").Value = 2.2
$ws.Range("I6").Value = 2.8
$ws.Range("J6").Value = 2.63
$ws.Range("K6").Value = 2.4
$ws.Range("L6").Value = 3.1
$ws.Range("Q6").Value = 1.44
$ws.Range("R6").Value = 2.63
$ws.Range("U6").Value = 1.44
$ws.Range("V6").Value = 2.63
$ws.Range("W6").Value = 13
$ws.Range("AA6").Value = 15
$ws.Range("AC6").Value = 19
$ws.Range("AL6").Value = 21
$ws.Range("AN6").Value = 4.75
$ws.Range("AO6").Value = 11
$ws.Range("AY6").Value = 19

# --- 6) Row 7 after the shift (Grasshoppers - Lugano) ---------------------
$ws.Range("U7").Value = 1.63

# --- 7) Row 8 after the shift (Servette - Luzern) -------------------------
$ws.Range("G8").Value = 1.8
$ws.Range("H8").Value = 3.75
$ws.Range("I8").Value = 4.1
$ws.Range("J8").Value = 2.38
$ws.Range("U8").Value = 1.5
$ws.Range("V8").Value = 2.37
$ws.Range("AC8").Value = 17
$ws.Range("AG8").Value = 101
$ws.Range("AH8").Value = 17
$ws.Range("AL8").Value = 29
$ws.Range("AN8").Value = 4
$ws.Range("AS8").Value = 81
